$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

$ws.Range("C2").Value = "aaa"
$ws.Range("C3").Value = "bbb"
$ws.Range("C4").Value = "ccc"
$ws.Range("C5").Value = "ddd"
$ws.Range("C6").Value = "eee"
$ws.Range("C7").Value = "fff"

$ws.Range("D26").Select()
